$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose Hungarian translation (column B) was empty:
# row 96  -> "Used in %(count)s unit" (no translation present)
# row 106 -> "Contains %(count)s group" (no translation present)
# Deleting row 106 first keeps row 96's index valid.
$ws.Rows("106").Delete()
$ws.Rows("96").Delete()
